$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 134983.5
$ws.Range("J17").Value = 137019.25
$ws.Range("L17").Value = 411057.75
$ws.Range("N17").Value = -411393.75

$ws.Range("H112").Value = 1783.1765
$ws.Range("I112").Value = 1100
$ws.Range("J112").Value = 1825.875
$ws.Range("K112").Value = 3300
$ws.Range("L112").Value = 5477.625
$ws.Range("M112").Value = -2192
$ws.Range("N112").Value = -7693.625

$ws.Range("H137").Value = 3475.1091
$ws.Range("I137").Value = 2632.9512
$ws.Range("J137").Value = 5941.4287
$ws.Range("K137").Value = 7898.8536
$ws.Range("L137").Value = 17824.2861
$ws.Range("M137").Value = -5348.8536
$ws.Range("N137").Value = -22924.2861

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 974.4545000000001
$ws.Range("I61").Value = 759.1579
$ws.Range("J61").Value = 2338
$ws.Range("K61").Value = 759.1579
$ws.Range("L61").Value = 2338
$ws.Range("M61").Value = -547.1579
$ws.Range("N61").Value = -2762

$ws.Range("H74").Value = 1101.1136
$ws.Range("I74").Value = 595.7692
$ws.Range("K74").Value = 595.7692
$ws.Range("M74").Value = 278.2308

$ws.Range("H77").Value = 1101.1136
$ws.Range("I77").Value = 595.7692
$ws.Range("K77").Value = 2978.846
$ws.Range("M77").Value = 1389.154

$ws.Range("H136").Value = 974.4545000000001
$ws.Range("I136").Value = 759.1579
$ws.Range("J136").Value = 2338
$ws.Range("K136").Value = 2277.4737
$ws.Range("L136").Value = 7014
$ws.Range("M136").Value = 272.5263
$ws.Range("N136").Value = -12114

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 80780
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 80780
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 80780
$ws.Range("M57").Value = ""
$ws.Range("N57").Value = -82220

$ws.Range("H59").Value = 53326.668
$ws.Range("J59").Value = 53326.668
$ws.Range("L59").Value = 53326.668
$ws.Range("N59").Value = -55020.668

$ws.Range("H60").Value = 69500
$ws.Range("J60").Value = 69500
$ws.Range("L60").Value = 69500
$ws.Range("N60").Value = -70698

$ws.Range("H86").Value = 3726.348
$ws.Range("I86").Value = 3575.5
$ws.Range("J86").Value = 4071.1428
$ws.Range("K86").Value = 3575.5
$ws.Range("L86").Value = 4071.1428
$ws.Range("M86").Value = -2452.5
$ws.Range("N86").Value = -6317.1428

$ws.Range("H89").Value = 3726.348
$ws.Range("I89").Value = 3575.5
$ws.Range("J89").Value = 4071.1428
$ws.Range("K89").Value = 17877.5
$ws.Range("L89").Value = 20355.714
$ws.Range("M89").Value = -12261.5
$ws.Range("N89").Value = -31587.714

$ws.Range("H134").Value = 900.4167
$ws.Range("I134").Value = 765.069
$ws.Range("J134").Value = 1461.1428
$ws.Range("K134").Value = 2295.207
$ws.Range("L134").Value = 4383.428400000001
$ws.Range("M134").Value = 239.7930000000001
$ws.Range("N134").Value = -9453.428400000001

$ws.Range("H136").Value = 80780
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 80780
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 80780
$ws.Range("M136").Value = ""
$ws.Range("N136").Value = -90980

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1174.2142
$ws.Range("I22").Value = 1348.8889
$ws.Range("J22").Value = 859.8
$ws.Range("K22").Value = 1348.8889
$ws.Range("L22").Value = 859.8
$ws.Range("M22").Value = -998.8888999999999
$ws.Range("N22").Value = -1559.8

$ws.Range("H31").Value = 22107.334
$ws.Range("I31").Value = 29337.162
$ws.Range("J31").Value = 2999.9285
$ws.Range("K31").Value = 29337.162
$ws.Range("L31").Value = 2999.9285
$ws.Range("M31").Value = -29042.162
$ws.Range("N31").Value = -3589.9285

$ws.Range("H34").Value = 22107.334
$ws.Range("I34").Value = 29337.162
$ws.Range("J34").Value = 2999.9285
$ws.Range("K34").Value = 29337.162
$ws.Range("L34").Value = 2999.9285
$ws.Range("M34").Value = -29135.162
$ws.Range("N34").Value = -3403.9285

$ws.Range("H58").Value = 2154.32
$ws.Range("I58").Value = 2084.5217
$ws.Range("J58").Value = 2957
$ws.Range("K58").Value = 2084.5217
$ws.Range("L58").Value = 2957
$ws.Range("M58").Value = -1881.5217
$ws.Range("N58").Value = -3363

$ws.Range("H127").Value = 35309.09
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 35309.09
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 35309.09
$ws.Range("M127").Value = ""
$ws.Range("N127").Value = -45229.09

$ws.Range("H132").Value = 1150.375
$ws.Range("I132").Value = 940.6
$ws.Range("J132").Value = 2008.5454
$ws.Range("K132").Value = 2821.8
$ws.Range("L132").Value = 6025.6362
$ws.Range("M132").Value = -291.8000000000002
$ws.Range("N132").Value = -11085.6362

$ws.Range("H134").Value = 1413.7778
$ws.Range("I134").Value = 1420.0588
$ws.Range("J134").Value = 1307
$ws.Range("K134").Value = 4260.1764
$ws.Range("L134").Value = 3921
$ws.Range("M134").Value = -1725.1764
$ws.Range("N134").Value = -8991

$ws.Range("H136").Value = 2154.32
$ws.Range("I136").Value = 2084.5217
$ws.Range("J136").Value = 2957
$ws.Range("K136").Value = 6253.5651
$ws.Range("L136").Value = 8871
$ws.Range("M136").Value = -3703.5651
$ws.Range("N136").Value = -13971

$ws.Range("H137").Value = 73950
$ws.Range("J137").Value = 73950
$ws.Range("L137").Value = 73950
$ws.Range("N137").Value = -84150

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H123").Value = 3271.625
$ws.Range("I123").Value = 1328.25
$ws.Range("J123").Value = 3919.4167
$ws.Range("K123").Value = 3984.75
$ws.Range("L123").Value = 11758.2501
$ws.Range("M123").Value = -1534.75
$ws.Range("N123").Value = -16658.2501

$ws.Range("H129").Value = 4274717.5
$ws.Range("I129").Value = 682.5
$ws.Range("J129").Value = 5051814.5
$ws.Range("K129").Value = 2047.5
$ws.Range("L129").Value = 15155443.5
$ws.Range("M129").Value = 2952.5
$ws.Range("N129").Value = -15165443.5

$ws.Range("H130").Value = 3217.2666
$ws.Range("J130").Value = 3302.0715
$ws.Range("L130").Value = 9906.2145
$ws.Range("N130").Value = -19946.2145

$ws.Range("H131").Value = 812.987
$ws.Range("I131").Value = 519.875
$ws.Range("J131").Value = 889.86884
$ws.Range("K131").Value = 1559.625
$ws.Range("L131").Value = 2669.60652
$ws.Range("M131").Value = 3480.375
$ws.Range("N131").Value = -12749.60652

$ws.Range("H133").Value = 4385.423
$ws.Range("I133").Value = 5322.857
$ws.Range("J133").Value = 4040.0527
$ws.Range("K133").Value = 15968.571
$ws.Range("L133").Value = 12120.1581
$ws.Range("M133").Value = -10908.571
$ws.Range("N133").Value = -22240.1581

$ws.Range("H134").Value = 2619.2666
$ws.Range("I134").Value = 977.8182
$ws.Range("K134").Value = 2933.4546
$ws.Range("M134").Value = 2136.5454

$ws.Range("H136").Value = 2989.4119
$ws.Range("I136").Value = 1414
$ws.Range("J136").Value = 3645.8333
$ws.Range("K136").Value = 4242
$ws.Range("L136").Value = 10937.4999
$ws.Range("M136").Value = 858
$ws.Range("N136").Value = -21137.4999

$ws.Range("H137").Value = 3388.2144
$ws.Range("I137").Value = 1204.875
$ws.Range("J137").Value = 4261.55
$ws.Range("K137").Value = 3614.625
$ws.Range("L137").Value = 12784.65
$ws.Range("M137").Value = 1485.375
$ws.Range("N137").Value = -22984.65

$ws.Range("H138").Value = 2163.24
$ws.Range("I138").Value = 938.5
$ws.Range("J138").Value = 2550
$ws.Range("K138").Value = 2815.5
$ws.Range("L138").Value = 7650
$ws.Range("M138").Value = 2324.5
$ws.Range("N138").Value = -17930

$ws.Range("H139").Value = 1432.6
$ws.Range("I139").Value = 798.9
$ws.Range("J139").Value = 2700
$ws.Range("K139").Value = 2396.7
$ws.Range("L139").Value = 8100
$ws.Range("M139").Value = 2743.3
$ws.Range("N139").Value = -18380

$ws.Range("H140").Value = 1728.1666
$ws.Range("I140").Value = 1330.3214
$ws.Range("J140").Value = 3120.625
$ws.Range("K140").Value = 3990.9642
$ws.Range("L140").Value = 9361.875
$ws.Range("M140").Value = 1189.0358
$ws.Range("N140").Value = -19721.875

$ws.Range("H141").Value = 4671.875
$ws.Range("I141").Value = 3175
$ws.Range("J141").Value = 7166.6665
$ws.Range("K141").Value = 9525
$ws.Range("L141").Value = 21499.9995
$ws.Range("M141").Value = -4345
$ws.Range("N141").Value = -31859.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2215.9167
$ws.Range("I122").Value = 1659.1
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 4977.299999999999
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -2527.299999999999
$ws.Range("N122").Value = -19900

$ws.Range("H126").Value = 1217.909
$ws.Range("I126").Value = 1002.36365
$ws.Range("J126").Value = 1433.4546
$ws.Range("K126").Value = 3007.09095
$ws.Range("L126").Value = 4300.3638
$ws.Range("M126").Value = -537.0909499999998
$ws.Range("N126").Value = -9240.363799999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2422.818
$ws.Range("I16").Value = 1183.4445
$ws.Range("J16").Value = 8000
$ws.Range("K16").Value = 1183.4445
$ws.Range("L16").Value = 8000
$ws.Range("M16").Value = -1013.4445
$ws.Range("N16").Value = -8340

$ws.Range("H132").Value = 2033.1111
$ws.Range("I132").Value = 1740.7391
$ws.Range("J132").Value = 3714.25
$ws.Range("K132").Value = 5222.2173
$ws.Range("L132").Value = 11142.75
$ws.Range("M132").Value = -2692.2173
$ws.Range("N132").Value = -16202.75

$ws.Range("H136").Value = 2789.3125
$ws.Range("I136").Value = 2125.8708
$ws.Range("J136").Value = 3999.1177
$ws.Range("K136").Value = 6377.6124
$ws.Range("L136").Value = 11997.3531
$ws.Range("M136").Value = -3827.6124
$ws.Range("N136").Value = -17097.3531

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 375470.84
$ws.Range("I122").Value = 422050.53
$ws.Range("J122").Value = 2833.3333
$ws.Range("K122").Value = 1266151.59
$ws.Range("L122").Value = 8499.999899999999
$ws.Range("M122").Value = -1263701.59
$ws.Range("N122").Value = -13399.9999

$ws.Range("H132").Value = 451.82257
$ws.Range("I132").Value = 362.55356
$ws.Range("J132").Value = 1285
$ws.Range("K132").Value = 1087.66068
$ws.Range("L132").Value = 3855
$ws.Range("M132").Value = 1442.33932
$ws.Range("N132").Value = -8915

$ws.Range("H136").Value = 401.3256
$ws.Range("I136").Value = 299.94446
$ws.Range("J136").Value = 922.7143
$ws.Range("K136").Value = 899.83338
$ws.Range("L136").Value = 2768.1429
$ws.Range("M136").Value = 1650.16662
$ws.Range("N136").Value = -7868.1429

Write-Host "Applied all updates"
